# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-27 14:15:32
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet contains
# comma-separated lists of recorders. Several rows had the order of the names in
# that list swapped (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# Apply the same exact-value remapping to every matching cell in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact old-value -> new-value remapping observed for the "Recorded By" column.
$mapping = @{
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

# Determine the used range of the sheet so we cover every data row in column G.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($mapping.ContainsKey($val)) {
        $cell.Value2 = $mapping[$val]
    }
}
